$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2186.5557
$ws.Range("I58").Value = 316.66666
$ws.Range("J58").Value = 2560.5334
$ws.Range("K58").Value = 949.9999799999999
$ws.Range("L58").Value = 7681.600199999999
$ws.Range("M58").Value = -799.9999799999999
$ws.Range("N58").Value = -7981.600199999999
$ws.Range("H103").Value = 1574.8125
$ws.Range("I103").Value = 1171.4286
$ws.Range("J103").Value = 1888.5555
$ws.Range("K103").Value = 3514.2858
$ws.Range("L103").Value = 5665.666499999999
$ws.Range("M103").Value = -2928.2858
$ws.Range("N103").Value = -6837.666499999999
$ws.Range("H125").Value = 1145.3334
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1145.3334
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 10308.0006
$ws.Range("N125").Value = -15228.0006
$ws.Range("M125").ClearContents()
$ws.Range("H132").Value = 6445.9287
$ws.Range("I132").Value = 5412.654
$ws.Range("K132").Value = 16237.962
$ws.Range("M132").Value = -13707.962

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2578.879
$ws.Range("I61").Value = 1484.0526
$ws.Range("J61").Value = 4064.7144
$ws.Range("K61").Value = 1484.0526
$ws.Range("L61").Value = 4064.7144
$ws.Range("M61").Value = -1272.0526
$ws.Range("N61").Value = -4488.7144
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 1654.4762
$ws.Range("I74").Value = 1560.2424
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1560.2424
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -686.2424000000001
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1654.4762
$ws.Range("I77").Value = 1560.2424
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 7801.212
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -3433.212
$ws.Range("N77").Value = -18736
$ws.Range("H97").Value = 1007.0417
$ws.Range("I97").Value = 771.55554
$ws.Range("K97").Value = 771.55554
$ws.Range("M97").Value = -275.55554
$ws.Range("H110").Value = 1686.8077
$ws.Range("I110").Value = 1502.7142
$ws.Range("J110").Value = 2460
$ws.Range("K110").Value = 1502.7142
$ws.Range("L110").Value = 2460
$ws.Range("M110").Value = 542.2858000000001
$ws.Range("N110").Value = -6550
$ws.Range("H132").Value = 4361.9316
$ws.Range("I132").Value = 1707.4286
$ws.Range("J132").Value = 9007.3125
$ws.Range("K132").Value = 5122.2858
$ws.Range("L132").Value = 27021.9375
$ws.Range("M132").Value = -2592.2858
$ws.Range("N132").Value = -32081.9375
$ws.Range("H136").Value = 2578.879
$ws.Range("I136").Value = 1484.0526
$ws.Range("J136").Value = 4064.7144
$ws.Range("K136").Value = 4452.1578
$ws.Range("L136").Value = 12194.1432
$ws.Range("M136").Value = -1902.1578
$ws.Range("N136").Value = -17294.1432

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4569.6035
$ws.Range("I134").Value = 2196.5386
$ws.Range("J134").Value = 6497.7188
$ws.Range("K134").Value = 6589.6158
$ws.Range("L134").Value = 19493.1564
$ws.Range("M134").Value = -4054.6158
$ws.Range("N134").Value = -24563.1564

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H122").Value = 71430456
$ws.Range("I122").Value = 111111980
$ws.Range("J122").Value = 3699.8
$ws.Range("K122").Value = 333335940
$ws.Range("L122").Value = 11099.4
$ws.Range("M122").Value = -333333490
$ws.Range("N122").Value = -15999.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2886.3125
$ws.Range("I80").Value = 950.25
$ws.Range("J80").Value = 3531.6667
$ws.Range("K80").Value = 2850.75
$ws.Range("L80").Value = 10595.0001
$ws.Range("M80").Value = -1914.75
$ws.Range("N80").Value = -12467.0001
$ws.Range("H83").Value = 2886.3125
$ws.Range("I83").Value = 950.25
$ws.Range("J83").Value = 3531.6667
$ws.Range("K83").Value = 8552.25
$ws.Range("L83").Value = 31785.0003
$ws.Range("M83").Value = -3872.25
$ws.Range("N83").Value = -41145.0003
$ws.Range("H113").Value = 728.129
$ws.Range("I113").Value = 572.625
$ws.Range("K113").Value = 1717.875
$ws.Range("M113").Value = 452.125
$ws.Range("H115").Value = 2728
$ws.Range("I115").Value = 2141.5
$ws.Range("J115").Value = 3607.75
$ws.Range("K115").Value = 6424.5
$ws.Range("L115").Value = 10823.25
$ws.Range("M115").Value = -5249.5
$ws.Range("N115").Value = -13173.25
$ws.Range("H134").Value = 2809.1538
$ws.Range("I134").Value = 1668.7778
$ws.Range("J134").Value = 5375
$ws.Range("K134").Value = 5006.3334
$ws.Range("L134").Value = 16125
$ws.Range("M134").Value = 63.66659999999956
$ws.Range("N134").Value = -26265

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3108261.5
$ws.Range("I102").Value = 5496194
$ws.Range("J102").Value = 3949.2
$ws.Range("K102").Value = 5496194
$ws.Range("L102").Value = 3949.2
$ws.Range("M102").Value = -5494572
$ws.Range("N102").Value = -7193.2
$ws.Range("H132").Value = 2818.4517
$ws.Range("I132").Value = 2581
$ws.Range("J132").Value = 2915.5908
$ws.Range("K132").Value = 7743
$ws.Range("L132").Value = 8746.7724
$ws.Range("M132").Value = -5213
$ws.Range("N132").Value = -13806.7724

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2142.8572
$ws.Range("I7").Value = 2125
$ws.Range("J7").Value = 2166.6667
$ws.Range("K7").Value = 2125
$ws.Range("L7").Value = 2166.6667
$ws.Range("M7").Value = -2013
$ws.Range("N7").Value = -2390.6667
$ws.Range("H100").Value = 3629.7778
$ws.Range("I100").Value = 3490
$ws.Range("J100").Value = 3683.5386
$ws.Range("K100").Value = 3490
$ws.Range("L100").Value = 3683.5386
$ws.Range("M100").Value = -2949
$ws.Range("N100").Value = -4765.5386
$ws.Range("H126").Value = 2142.8572
$ws.Range("I126").Value = 2125
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 6375
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -3905
$ws.Range("N126").Value = -11440.0001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3042.5652
$ws.Range("I126").Value = 2978.2666
$ws.Range("J126").Value = 3163.125
$ws.Range("K126").Value = 8934.799800000001
$ws.Range("L126").Value = 9489.375
$ws.Range("M126").Value = -6464.799800000001
$ws.Range("N126").Value = -14429.375
$ws.Range("H132").Value = 2625.081
$ws.Range("I132").Value = 2630.1875
$ws.Range("K132").Value = 7890.5625
$ws.Range("M132").Value = -5360.5625
